# Actualización automática 2025-11-26 16:30:11
# A new client/asesor record ("GUERRERO GARCIA OLIMPIA ANNABELLE") is
# inserted alphabetically (between "GRANJA VANEGAS MARCELA" and
# "JAIME COELLO ALBERTO FERNANDO") at row 27 of both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets. This pushes every
# following row down by one and grows the used range by one row.
# The new row carries zero values in every numeric column, and the
# trailing "X de 55" summary counters on "VENTAS POR GRUPO" need their
# denominator bumped to "X de 56" to reflect the extra row (the
# numeric sum row on "VENTA MENSUAL" keeps the same totals, Excel
# recomputes that automatically because it is the SUM of the column).

$wb = $excel.ActiveWorkbook

$newName = "GUERRERO GARCIA OLIMPIA ANNABELLE"
$asesor  = "OFICINA-CATAECSA"

# ---- Sheet 1: "VENTAS POR GRUPO" (columns A:R) ----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new blank row at 27, shifting rows 27:57 down to 28:58.
$ws1.Rows.Item(27).Insert()

# Populate the newly-inserted row with the new record (all metrics 0).
$ws1.Range("A27").Value = $asesor
$ws1.Range("B27").Value = $newName
$ws1.Range("C27:R27").Value = 0

# The trailing counter row (now row 58) reads like "0 de 55" for every
# metric column; bump the denominator to reflect the new row count (56).
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(58, $col)
    $text = $cell.Value()
    $cell.Value = ($text -replace "de 55", "de 56")
}

# ---- Sheet 2: "VENTA MENSUAL" (columns A:G) ----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Insert a new blank row at 27, shifting rows 27:57 down to 28:58.
$ws2.Rows.Item(27).Insert()

# Populate the newly-inserted row with the new record (all metrics 0).
$ws2.Range("A27").Value = $asesor
$ws2.Range("B27").Value = $newName
$ws2.Range("C27:G27").Value = 0
